$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text even when the new value looks numeric,
# then restore the default (Normal) style so no stray formatting is introduced.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.476.30'
$ws.Range("E2").Value = '  +2.32%  '
$ws.Range("D3").Value = '2.554.41'
$ws.Range("E3").Value = '  +5.29%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '571.89'
$ws.Range("E5").Value = '  +2.69%  '
$ws.Range("D6").Value = '149.86'
$ws.Range("E6").Value = '  +8.21%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +0.67%  '
$ws.Range("D9").Value = '2.553.10'
$ws.Range("E9").Value = '  +5.46%  '
$ws.Range("E10").Value = '  +2.35%  '
$ws.Range("E11").Value = '  -0.34%  '
$ws.Range("E12").Value = '  +1.36%  '
$ws.Range("D13").Value = '0.359'
$ws.Range("E13").Value = '  +3.62%  '
$ws.Range("D14").Value = '28.08'
$ws.Range("E14").Value = '  +9.21%  '
$ws.Range("D15").Value = '3.011.64'
$ws.Range("E15").Value = '  +5.50%  '
$ws.Range("D16").Value = '63.448.61'
$ws.Range("E16").Value = '  +2.36%  '
$ws.Range("D17").Value = '0.0000144'
$ws.Range("E17").Value = '  +2.81%  '
$ws.Range("D18").Value = '2.569.34'
$ws.Range("E18").Value = '  +6.01%  '
$ws.Range("D19").Value = '11.62'
$ws.Range("E19").Value = '  +4.64%  '
$ws.Range("D20").Value = '341.60'
$ws.Range("E20").Value = '  -0.94%  '
$ws.Range("E21").Value = '  +3.30%  '
$ws.Range("E22").Value = '  +1.11%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '66.21'
$ws.Range("E24").Value = '  +1.86%  '
$ws.Range("E25").Value = '  -1.48%  '
$ws.Range("E26").Value = '  +4.09%  '
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("D28").Value = '8.44'
$ws.Range("E28").Value = '  +2.10%  '
$ws.Range("E29").Value = '  +7.72%  '
$ws.Range("D30").Value = '7.23'
$ws.Range("E30").Value = '  +14.53%  '
$ws.Range("D31").Value = '0.0₃0838'
$ws.Range("E31").Value = '  +6.17%  '
$ws.Range("E32").Value = '  +3.58%  '
$ws.Range("D33").Value = '177.61'
$ws.Range("E33").Value = '  +3.93%  '
$ws.Range("E34").Value = '  +9.98%  '
$ws.Range("D35").Value = '414.52'
$ws.Range("E35").Value = '  +10.83%  '
$ws.Range("D36").Value = '0.403'
$ws.Range("E36").Value = '  +2.08%  '
$ws.Range("D37").Value = '19.09'
$ws.Range("E37").Value = '  +3.01%  '
$ws.Range("E38").Value = '  -0.44%  '
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("E40").Value = '  +3.91%  '
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("D42").Value = '39.99'
$ws.Range("E42").Value = '  +2.32%  '
$ws.Range("D43").Value = '155.20'
$ws.Range("E43").Value = '  +6.70%  '
$ws.Range("D44").Value = '3.80'
$ws.Range("E44").Value = '  +3.74%  '
$ws.Range("D45").Value = '21.09'
$ws.Range("E45").Value = '  +2.04%  '
$ws.Range("E46").Value = '  +4.03%  '
$ws.Range("D47").Value = '0.0533'
$ws.Range("E47").Value = '  +2.87%  '
$ws.Range("D48").Value = '0.0966'
$ws.Range("E48").Value = '  +1.05%  '
$ws.Range("E49").Value = '  +5.59%  '
$ws.Range("D50").Value = '18.75'
$ws.Range("E50").Value = '  +4.20%  '
$ws.Range("D51").Value = '1.86'
$ws.Range("E51").Value = '  +8.31%  '

$ws.Range("D2:D51").Style = "Normal"
